# Applies the diff: adds header label "Unnamed: 0" to A1 (with the bold/
# bordered header style), strips that same style from A2:A14 (they become
# plain/unstyled), and appends two new summary rows (15 "mean", 16 "std")
# with per-column statistics across B:Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: new header label "Unnamed: 0" ------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"

# Give A1 the same style as the other header cells (bold, centered, boxed)
# by copying formatting from an already-styled header cell (A2 still has it
# at this point) and pasting only the formats onto A1.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- A2:A14: drop the bold/boxed header style -> back to default "Normal"
$ws.Range("A2:A14").Style = "Normal"

# --- Row 15 ("mean") and Row 16 ("std") ------------------------------------
$ws.Range("A15").Value = "mean"
$ws.Range("A16").Value = "std"

$meanVals = @(0.4932613673042583,0.5382312831640398,0.4932613673042583,0.5371201720529286,0.4921502561931472,0.5382312831640398,0.4932613673042583,0.5382312831640398,0.4921502561931472,0.5371201720529286,0.4932613673042583,0.5382312831640398,0.4932613673042583,0.5382312831640398,0.4921502561931472,0.5382312831640398)
$stdVals  = @(0.3627495052232896,0.3487818395416044,0.3627495052232896,0.3506857010387602,0.3644262175929625,0.3487818395416044,0.3627495052232896,0.3487818395416044,0.3644262175929625,0.3506857010387602,0.3627495052232896,0.3487818395416044,0.3627495052232896,0.3487818395416044,0.3644262175929625,0.3487818395416044)

for ($i = 0; $i -lt $meanVals.Length; $i++) {
    $col = $i + 2  # B = 2 .. Q = 17
    $ws.Cells.Item(15, $col).Value = $meanVals[$i]
    $ws.Cells.Item(16, $col).Value = $stdVals[$i]
}
